$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (rows are 1-indexed in the Word COM model)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "1509"
$t.Cell(5,1).Range.Text  = "0.00001"
$t.Cell(6,1).Range.Text  = "0.00073"
$t.Cell(7,1).Range.Text  = "0.00012"
$t.Cell(8,1).Range.Text  = "0.00006"
$t.Cell(9,1).Range.Text  = "0.00015"
$t.Cell(10,1).Range.Text = "0.00016"
$t.Cell(11,1).Range.Text = "0.00017"
$t.Cell(12,1).Range.Text = "0.18515"

# Collapse the multi-run tab-separated cells down to a single value
$t.Cell(44,1).Range.Text = "99.91"
$t.Cell(45,1).Range.Text = "0.19"
$t.Cell(46,1).Range.Text = "199"
